# The row for "even_MAG-GUT49526.fa" (row 4) was removed from the
# output, shifting all subsequent rows up by one. Delete that entire
# row so the remaining data shifts up and the sheet's used range
# shrinks from A1:L9 to A1:L8.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("4:4").Delete()
